$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.508.53'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '2.060.36'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.13'
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.34'
$ws.Range('E8').Value = '  -1.69%  '
$ws.Range('E9').Value = '  -0.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0790'
$ws.Range('E10').Value = '  +0.70%  '
$ws.Range('E11').Value = '  -2.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.01'
$ws.Range('E12').Value = '  +1.76%  '
$ws.Range('D13').Value = '2.365.96'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.95'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.760'
$ws.Range('E15').Value = '  -2.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.32'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').Value = '2.060.36'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').Value = '37.414.51'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.13'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.95'
$ws.Range('E20').Value = '  -2.01%  '
$ws.Range('D21').Value = '0.0₃0831'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.89'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.32'
$ws.Range('E25').Value = '  -3.19%  '
$ws.Range('E26').Value = '  +1.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.85'
$ws.Range('E27').Value = '  -2.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.130'
$ws.Range('E28').Value = '  -5.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.18'
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('E30').Value = '  -2.25%  '
$ws.Range('E31').Value = '  -1.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.58'
$ws.Range('E32').Value = '  -3.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0623'
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.48'
$ws.Range('E35').Value = '  +1.21%  '
$ws.Range('E36').Value = '  +1.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.29'
$ws.Range('E37').Value = '  -3.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.26'
$ws.Range('E39').Value = '  -3.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0222'
$ws.Range('E40').Value = '  -4.51%  '
$ws.Range('D41').Value = '1.507.34'
$ws.Range('E41').Value = '  +4.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.27'
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.91'
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.65'
$ws.Range('E44').Value = '  -2.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0950'
$ws.Range('E45').Value = '  -2.78%  '
$ws.Range('E46').Value = '  +2.09%  '
$ws.Range('E47').Value = '  -2.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.97'
$ws.Range('E48').Value = '  -2.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.16'
$ws.Range('E49').Value = '  -2.34%  '
$ws.Range('E50').Value = '  -1.65%  '
$ws.Range('D51').Value = '2.251.12'
$ws.Range('E51').Value = '  -0.58%  '
